$d = $word.ActiveDocument

# --- Change 1: paragraph 5 (TU NAPSAC O EUKLIDESIE -> Euclid description) ---
$p5 = $d.Paragraphs(5)
$xml5 = '<w:p w14:paraId="3E968CB5" w14:textId="1C3B20AC" w:rsidR="00924B68" w:rsidRDefault="00924B68" w:rsidP="00F55889"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">Pierwszą implementacją w programie jest wyznaczanie NWD za pomocą algorytmu Euklidesa. Aby zobrazować działanie algorytmu, załóżmy że należy wyznaczyć NWD z liczb a oraz b. Na początku wykonywane jest dzielenie z resztą liczby a przez liczbę b. Jest to realizowane za pomocą instrukcji div w języku assembler. </w:t></w:r><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>Gdy r</w:t></w:r><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>eszta z dzielenia</w:t></w:r><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">, która jest </w:t></w:r><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>umieszcz</w:t></w:r><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>ona</w:t></w:r><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> w rejestrze %</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>edx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> wynosi 0 to największym wspólnym dzielnikiem jest liczba b. W przypadku gdy reszta jest różna od zera to następuje przypisanie liczbie a wartości liczby b. Następnie liczbie b jest przypisywana wartość reszty. Następnie ponownie jest realizowane dzielenie liczby a przez b, aż reszta nie będzie równa zero. Cała operacja jest realizowana w pętli </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>while</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">, poprzez użycie instrukcji </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>cmp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>. Implementacja algorytmu znajduje się w pliku Euklides.cpp. Dodatkowo znajduje się tam również funkcja wyświetlająca menu oraz pobierająca dane od użytkownika.</w:t></w:r></w:p>'
$p5.Range.InsertXML($xml5)
Write-Host "After change 1: count=$($d.Paragraphs.Count)"

# --- Change 2: paragraph 6 (text fix + bookmark removal) ---
$p6 = $d.Paragraphs(6)
$xml6 = '<w:p w14:paraId="4F575DDD" w14:textId="64967862" w:rsidR="008C18B0" w:rsidRPr="00F55889" w:rsidRDefault="00F55889" w:rsidP="00F55889"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">Kolejnym algorytmem, który został zaimplementowany jest naiwny algorytm wyznaczania pierwszości liczby. Algorytm polega na </w:t></w:r><w:r w:rsidR="004F7993"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">próbnym </w:t></w:r><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">dzieleniu sprawdzanej liczby </w:t></w:r><w:r w:rsidRPr="00F55889"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>a</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>przez liczby z</w:t></w:r><w:r w:rsidR="00144759"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t> </w:t></w:r><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">zakresu od </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>2</w:t></w:r><w:r w:rsidRPr="00F55889"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">do </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F55889"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>sqrt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00F55889"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>(a)</w:t></w:r><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidR="004F7993"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> Przy każdej takiej operacji badana jest reszta z dzielenia</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>. Jeśli reszta podczas które</w:t></w:r><w:r w:rsidR="00C97D11"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">jś operacji dzielenia </w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>będzie wynosiła 0, to liczba</w:t></w:r><w:r w:rsidR="00C97D11"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C97D11" w:rsidRPr="00C97D11"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>a</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> nie będzie liczbą pierwszą. Interesujące może się wydawać, dlaczego wystarczy sprawdzić liczby z podanego wyżej zakresu. Dzieje się tak, ponieważ jeśli liczba posiada czynnik większy od </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="008C18B0" w:rsidRPr="00AD323B"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>sqrt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="008C18B0" w:rsidRPr="00AD323B"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>(a)</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>, to drugi jego czynnik musi być mniejszy od pierwiastka</w:t></w:r><w:r w:rsidR="00C97D11"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> z </w:t></w:r><w:r w:rsidR="00C97D11" w:rsidRPr="00C97D11"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>a</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">, aby ich iloczyn musiał być równy </w:t></w:r><w:r w:rsidR="008C18B0" w:rsidRPr="00C97D11"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>a</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">. Zatem wystarczy podzielić liczbę </w:t></w:r><w:r w:rsidR="008C18B0" w:rsidRPr="00C97D11"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>a</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> przez liczby z danego </w:t></w:r><w:r w:rsidR="00C97D11"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>przedziału,</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> aby wykluczyć liczby złożone.  Tak jak zostało </w:t></w:r><w:r w:rsidR="00C97D11"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>wcześniej założone</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">, algorytm został napisany w języku asembler. Z powodu braku dobrej </w:t></w:r><w:r w:rsidR="00C97D11"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>znajomości</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> tego języka, pierwiastkowanie liczby </w:t></w:r><w:r w:rsidR="008C18B0" w:rsidRPr="00C97D11"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>a</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> zostało wykonane w języku </w:t></w:r><w:r w:rsidR="008C18B0" w:rsidRPr="00C97D11"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>C++.</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> Wynik tej operacji został jednak wykorzystany bezpośrednio we wstawionym fragmencie z kodem algorytmu. A</w:t></w:r><w:r w:rsidR="00C97D11"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>l</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">gorytm zwraca podmienioną liczbę </w:t></w:r><w:r w:rsidR="008C18B0" w:rsidRPr="00C97D11"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>a</w:t></w:r><w:r w:rsidR="00C97D11"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>, w zależności od tego czy jest pierwsza, czy nie</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">. Jeśli liczba </w:t></w:r><w:r w:rsidR="008C18B0" w:rsidRPr="00C97D11"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>a</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> była pierwsza, to zostanie zwrócone </w:t></w:r><w:r w:rsidR="008C18B0" w:rsidRPr="00C97D11"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>1</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">, a jeśli nie była liczbą pierwszą to zwróci </w:t></w:r><w:r w:rsidR="008C18B0" w:rsidRPr="00C97D11"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>0</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>. Zostało to wykorzy</w:t></w:r><w:r w:rsidR="00C97D11"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>st</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>ane w</w:t></w:r><w:r w:rsidR="00144759"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t> </w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">wywołaniu algorytmu. W pliku </w:t></w:r><w:r w:rsidR="00077713" w:rsidRPr="00077713"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>NaivePrime</w:t></w:r><w:r w:rsidR="008C18B0" w:rsidRPr="008C18B0"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>.cpp</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>,</w:t></w:r><w:r w:rsidR="00C97D11"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> w</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> którym znajduje się implementacja algorytmu, znajduje się również funkcja odpowiedzialna za wczytywanie liczby podanej obserwacji</w:t></w:r><w:r w:rsidR="00C97D11"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C97D11"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>W</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">ykorzystując zaimplementowany algorytm w funkcji </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="0048353D"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>c</w:t></w:r><w:r w:rsidR="00077713"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>heckPrime</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00077713"><w:rPr><w:i/><w:iCs/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve">() </w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> pokazuje</w:t></w:r><w:r w:rsidR="00C97D11"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> ona</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>, czy dana liczba jest liczbą pierwszą</w:t></w:r><w:r w:rsidR="004704BF"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> wyświetlając stosowną informację</w:t></w:r><w:r w:rsidR="008C18B0"><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>.</w:t></w:r></w:p>'
$p6.Range.InsertXML($xml6)
Write-Host "After change 2: count=$($d.Paragraphs.Count)"

# --- Change 3: merge paragraphs 7 & 8, add new text + bookmark ---
$p7 = $d.Paragraphs(7)
$markRange = $d.Range($p7.Range.End - 1, $p7.Range.End)
$markRange.Delete()
Write-Host "After merge: count=$($d.Paragraphs.Count)"
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$xml78 = '<w:p w14:paraId="13F1D9A3" w14:textId="1B844A8C" w:rsidR="00F55889" w:rsidRPr="00F55889" w:rsidRDefault="00F55889" w:rsidP="00F55889"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>W następnym etapie projektu mamy zamiar zaimplementować kolejne algorytmy. Będą  one bardziej złożone co powoduje że napisanie pełnych funkcji w języku assembler może być problematyczne. Realizowany jednak projekt pozwala na poznanie praktyczne assemblera, a więc będziemy dalej próbować implementować poszczególne funkcje przy jego użyciu.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$finalPara.Range.InsertXML($xml78)
Write-Host "After change 3: count=$($d.Paragraphs.Count)"
